$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Not Litigated")

# The "Ignored, Death" / Marilyn Monroe case currently sits in row 3.
# It needs to move down so it becomes the last row of this block (row 12),
# with every row currently between 4 and 12 shifting up by one.
# Equivalent to: delete row 3 (shifting 4..13 up to 3..12), then open a
# fresh blank row at 12 (pushing what is now row 12 back down to 13), and
# fill that blank row with the case's original data.

$ws.Rows("3:3").Delete()
$ws.Rows("12:12").Insert()

$ws.Range("A12").Value = "Ignored, Death"
$ws.Range("B12").Value = "California"
$ws.Range("C12").Value = "1950s-1962"
$ws.Range("D12").Value = "Fred Otash, Clients, Kennedys"
$ws.Range("E12").Value = "Marilyn Monroe"
$ws.Range("F12").Value = "Recording Audios, Fixing"
$ws.Range("L12").Value = "Film, Fixers, Investigators, Modeing"
$ws.Range("O12").Value = "Fixers"
$ws.Range("S12").Value = "https://linkedinvestigations.com/a-los-angeles-detective-caught-marilyn-monroe-s-death-on-tape/"

# The row delete/insert leaves the sheet's hyperlink collection stale (it
# still anchors the Marilyn Monroe link to the old S3 address, which now
# holds a different case). Rebuild both hyperlinks on the sheet so each
# lands on the correct cell, then restore the "Hyperlink" cell style
# (reusing the workbook's existing Hyperlink style instead of minting a
# new one) on both link cells.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("S12"), "https://linkedinvestigations.com/a-los-angeles-detective-caught-marilyn-monroe-s-death-on-tape/")
$ws.Hyperlinks.Add($ws.Range("S30"), "https://charleyproject.org/case/anthony-vivien-fox")
$ws.Range("S12").Style = "Hyperlink"
$ws.Range("S30").Style = "Hyperlink"

# Reflect the author's final selection on this sheet.
$ws.Range("C31").Select()
